$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 23

$ws.Cells.Item($row, 1).Value = "Record"
$ws.Cells.Item($row, 2).Value = "Balanço Geral"
$ws.Cells.Item($row, 3).Value = "Defesa Civil"
$ws.Cells.Item($row, 4).Value = "2025-04-01T13:06"
$ws.Cells.Item($row, 5).Value = "Neutro"
$ws.Cells.Item($row, 6).Value = "Defesa Civil de Campos realiza demolição parcial de prédio com risco de desabar. Repórter *ao vivo*. Vídeo com depoimento do secretário da defesa Civil, Alcemir Pascoutto. "
